$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B13's text (description of V7): append deprecation note
$ws.Range("B13").Value = "Without sliding window; deprecated, unsure how to build new model architecture"

# Add new row 15 for V9
$ws.Range("A15").Value = "V9"
$ws.Range("B15").Value = "Basically like model 6 but with abilty to read in config files"

# Update selection to match new active cell (A16)
$ws.Range("A16").Select()
